$wb = $excel.ActiveWorkbook

# Update "展览" sheet (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4941
$ws1.Range("F4").Value = 867

# Update "全部类型" sheet (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4941
$ws4.Range("F4").Value = 867
